# camparam.xlsx update: add MOG/MOG2 background-subtractor optical-flow
# method support ("implemented excellent working GMM code").
#
#   1. ofmethod default value changes from "farneback" to "mog2" (B2:C2).
#   2. The A2 cell comment documenting valid ofmethod values gains two
#      extra lines describing "mog" and "mog2".
#   3. Three new parameter rows are appended to the sheet: nhistory,
#      nmixtures and varThreshold (rows 26-28; row 25 is intentionally
#      left empty/untouched).
#   4. The sheet selection moves to A30 and the used range grows to
#      A1:C28 accordingly (handled automatically by the engine once the
#      new cells are written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. ofmethod default: farneback -> mog2
$ws.Range("B2").Value = "mog2"
$ws.Range("C2").Value = "mog2"

# 2. extend the A2 comment with the new mog / mog2 choices
$comment = $ws.Range("A2").Comment
$commentText = "hs: Horn-Schunck (opencv2 only)`nfarneback: Farneback method`nmog: background subtractor (opencv2 only)`nmog2: background subtractor"
[void]$comment.Text($commentText)

# 3. new GMM (background subtractor) parameters
$ws.Range("A26").Value = "nhistory"
$ws.Range("B26").Value = 100
$ws.Range("C26").Value = 100

$ws.Range("A27").Value = "nmixtures"
$ws.Range("B27").Value = 5
$ws.Range("C27").Value = 5

$ws.Range("A28").Value = "varThreshold"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 1

# 4. move the active selection to A30, matching the author's final cursor
[void]$ws.Range("A30").Select()

# Best-effort cosmetic tweak (bookViews@tabRatio 976 -> 983). Purely a
# view-state attribute with no bridged persistence path in this runtime;
# harmless if it is a no-op.
try {
  $win = $excel.ActiveWindow
  $win.TabRatio = 0.983
} catch {
}
